$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 233 ---
$ws.Range("D233").Value = 44595
$ws.Range("J233").Value = 75
$ws.Range("K233").Value = 20000
$ws.Range("L233").Value = 20000
$ws.Range("M233").Value = 20000
$ws.Range("O233").Value = "Región Metropolitana"
$ws.Range("P233").Value = 800

# --- Update existing row 234 ---
$ws.Range("D234").Value = 44595
$ws.Range("H234").Value = "Chilena(o)"
$ws.Range("J234").Value = 30
$ws.Range("K234").Value = 25000
$ws.Range("L234").Value = 25000
$ws.Range("M234").Value = 25000
$ws.Range("O234").Value = "Región Metropolitana"
$ws.Range("P234").Value = 1000

# --- Update existing row 235 ---
$ws.Range("D235").Value = 44335
$ws.Range("H235").Value = "Americana (o)"
$ws.Range("J235").Value = 16
$ws.Range("K235").Value = 33000
$ws.Range("L235").Value = 34000
$ws.Range("M235").Value = 33500
$ws.Range("P235").Value = 1340

# --- Add new row 236 (what used to be row 234's original data) ---
$ws.Range("A236").Value = 9
$ws.Range("B236").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C236").Value = "Metropolitana"
$ws.Range("D236").Value = 44552
$ws.Range("D236").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E236").Value = 13
$ws.Range("F236").Value = 100112021
$ws.Range("G236").Value = "Ají"
$ws.Range("H236").Value = "Americana (o)"
$ws.Range("I236").Value = "Primera"
$ws.Range("J236").Value = 7
$ws.Range("K236").Value = 25000
$ws.Range("L236").Value = 26000
$ws.Range("M236").Value = 25429
$ws.Range("N236").Value = "`$/caja 25 kilos"
$ws.Range("O236").Value = "Provincia de Limarí"
$ws.Range("P236").Value = 1017
$ws.Range("Q236").Value = 25
$ws.Range("R236").Value = "Hortaliza"

# --- Add new row 237 (what used to be row 235's original data) ---
$ws.Range("A237").Value = 9
$ws.Range("B237").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C237").Value = "Metropolitana"
$ws.Range("D237").Value = 44552
$ws.Range("D237").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E237").Value = 13
$ws.Range("F237").Value = 100112021
$ws.Range("G237").Value = "Ají"
$ws.Range("H237").Value = "Chilena(o)"
$ws.Range("I237").Value = "Primera"
$ws.Range("J237").Value = 7
$ws.Range("K237").Value = 59000
$ws.Range("L237").Value = 61000
$ws.Range("M237").Value = 59857
$ws.Range("N237").Value = "`$/caja 25 kilos"
$ws.Range("O237").Value = "Provincia de Huasco"
$ws.Range("P237").Value = 2394
$ws.Range("Q237").Value = 25
$ws.Range("R237").Value = "Hortaliza"
